$d = $word.ActiveDocument

function Get-ParaByText($doc, $text) {
    for ($i = 1; $i -le $doc.Paragraphs.Count; $i++) {
        $p = $doc.Paragraphs($i)
        $txt = $p.Range.Text.TrimEnd([char]13)
        if ($txt -eq $text) {
            return $p
        }
    }
    return $null
}

# ---------------------------------------------------------------------
# 1) "Name: TBD" -> "Name: " + "uroborus" (highlighted yellow)
# ---------------------------------------------------------------------
$d.Content.Find.Execute("Name: TBD", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Name: uroborus", 2)

$find = $d.Content.Find
$find.ClearFormatting()
$find.Replacement.ClearFormatting()
$find.Replacement.Highlight = $true
$find.Execute("uroborus", $false, $false, $false, $false, $false, $true, 1, `
              $false, "uroborus", 2)

# Reset the Find/Replace formatting state so it does not leak into later
# (unrelated) Find.Execute calls below.
$find.ClearFormatting()
$find.Replacement.ClearFormatting()

# ---------------------------------------------------------------------
# 2) "Player wins when the final tile is reached." -> strike-through
# ---------------------------------------------------------------------
$p = Get-ParaByText $d "Player wins when the final tile is reached."
$p.Range.Font.StrikeThrough = $true

# ---------------------------------------------------------------------
# 3) "Shop: sacrifice to gain benefit. " -> split into 3 runs with a
#    "_GoBack" bookmark between "(s)" and ". "; and
# 4) "Mystery: ..." -> drop the bookmark there, merge back to one run
# ---------------------------------------------------------------------
$d.Content.Find.Execute("Shop: sacrifice to gain benefit. ", $true, $false, `
                         $false, $false, $false, $true, 1, $false, `
                         "Shop: sacrifice to gain benefit(s). ", 2)

$p = Get-ParaByText $d "Shop: sacrifice to gain benefit(s). "
$shopStart = $p.Range.Start

# "Shop: sacrifice to gain benefit" = 31 chars, "(s)" = 3 chars, ". " = 2 chars
$split1 = $shopStart + 31
$split2 = $shopStart + 34

$bmRange1 = $d.Range($split1, $split1)
$d.Bookmarks.Add("TempSplit", $bmRange1)

$bmRange2 = $d.Range($split2, $split2)
$d.Bookmarks.Add("_GoBack", $bmRange2)

$d.Bookmarks("TempSplit").Delete()

# Merge the now-orphaned two runs in the "Mystery" paragraph back into one
# (the bookmark that used to split them has moved away, up into "Shop").
$d.Content.Find.Execute("Mystery: possibility for all, but not limited to, tile event.", `
                         $true, $false, $false, $false, $false, $true, 1, $false, `
                         "Mystery: possibility for all, but not limited to, tile event.", 2)
